$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new stock item ("OTRIVIN 0.1% ADULT NASAL DROPS 15 ML") is inserted as
# item #8, between the current item #7 (FLAGYL ..., row 13) and the current
# item #8 (ترمومتر زيئبق, row 14). Every item row below it shifts down by
# one row, the running-totals row shifts from row 18 to row 19 (with its
# price total bumped up by the new item's price), and the date/footer row
# shifts from row 19 to row 20.
#
# Every item data row (7..17) shares the exact same per-column style, so
# instead of Rows.Insert() (which recreates style/format records on this
# host), the row *contents* are copied downward one row at a time, walking
# from the bottom up so nothing gets clobbered before it is read. Columns
# that store text-look-alike numbers (H, L, N, P, Q) are re-entered with a
# leading apostrophe so they stay text cells, matching the source data.
# ---------------------------------------------------------------------------

function Copy-ItemRow($srcRow, $dstRow) {
    $ws.Range("C$dstRow").Value = "'" + [string]$ws.Range("C$srcRow").Value2
    $ws.Range("H$dstRow").Value = "'" + [string]$ws.Range("H$srcRow").Value2
    $ws.Range("L$dstRow").Value = "'" + [string]$ws.Range("L$srcRow").Value2
    $ws.Range("N$dstRow").Value = "'" + [string]$ws.Range("N$srcRow").Value2
    $ws.Range("P$dstRow").Value = "'" + [string]$ws.Range("P$srcRow").Value2
    $ws.Range("Q$dstRow").Value = "'" + [string]$ws.Range("Q$srcRow").Value2
}

# Footer row (old row 19) -> row 20
$ws.Range("A20").Value = "'" + [string]$ws.Range("A19").Value2
$ws.Range("G20").Value = "'" + [string]$ws.Range("G19").Value2
$ws.Range("K20").Value = "'" + [string]$ws.Range("K19").Value2

# Totals row (old row 18) -> row 19, with the grand total bumped for the
# newly added item's price (627.6 + 24.0 = 651.6)
$ws.Range("N19").Value = 651.6

# Item rows 17..14 shift down to 18..15
Copy-ItemRow 17 18
Copy-ItemRow 16 17
Copy-ItemRow 15 16
Copy-ItemRow 14 15

# Sequential item-number column (A) for rows 14..18 -> 8..12
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11
$ws.Range("A18").Value = 12

# New item #8 goes into row 14
$ws.Range("C14").Value = "OTRIVIN 0.1% ADULT NASAL DROPS 15 ML"
$ws.Range("H14").Value = "'3:0"
$ws.Range("L14").Value = "'1"
$ws.Range("N14").Value = "'24.00"
$ws.Range("P14").Value = "'24.0000"
$ws.Range("Q14").Value = "'1:0"

# The new row 18 (previously just a blank spacer row created by the shift)
# needs the same merged-cell layout as every other item row.
$ws.Range("A18:B18").Merge()
$ws.Range("C18:G18").Merge()
$ws.Range("H18:K18").Merge()
$ws.Range("L18:M18").Merge()
$ws.Range("N18:O18").Merge()

Write-Output "Row insert for OTRIVIN complete"
